$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 21:36"

# Update country rows: name swaps (re-sorted rank order) + refreshed stats
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 3329821
$ws.Cells.Item(4, 3).Value = 38035
$ws.Cells.Item(4, 4).Value = 1478542
$ws.Cells.Item(4, 5).Value = 1714105
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 503
$ws.Cells.Item(4, 8).Value = 137174

$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 850358
$ws.Cells.Item(6, 3).Value = 27755
$ws.Cells.Item(6, 4).Value = 536231
$ws.Cells.Item(6, 5).Value = 291440
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 543
$ws.Cells.Item(6, 8).Value = 22687

$ws.Cells.Item(9, 1).Value = "Chile"
$ws.Cells.Item(9, 2).Value = 312029
$ws.Cells.Item(9, 3).Value = 2755
$ws.Cells.Item(9, 4).Value = 281114
$ws.Cells.Item(9, 5).Value = 24034
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = 6881

$ws.Cells.Item(19, 1).Value = "Alemania"
$ws.Cells.Item(19, 2).Value = 199775
$ws.Cells.Item(19, 3).Value = 187
$ws.Cells.Item(19, 4).Value = 184500
$ws.Cells.Item(19, 5).Value = 6143
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 2
$ws.Cells.Item(19, 8).Value = 9132

$ws.Cells.Item(70, 1).Value = "Uzbekistan"
$ws.Cells.Item(70, 2).Value = 12513
$ws.Cells.Item(70, 3).Value = 486
$ws.Cells.Item(70, 4).Value = 7723
$ws.Cells.Item(70, 5).Value = 4733
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 3
$ws.Cells.Item(70, 8).Value = 57

$ws.Cells.Item(85, 1).Value = "Costa Rica"
$ws.Cells.Item(85, 2).Value = 7231
$ws.Cells.Item(85, 3).Value = 386
$ws.Cells.Item(85, 4).Value = 2220
$ws.Cells.Item(85, 5).Value = 4983
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 28

$ws.Cells.Item(86, 1).Value = "Bulgaria"
$ws.Cells.Item(86, 2).Value = 6964
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 3308
$ws.Cells.Item(86, 5).Value = 3389
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 267

$ws.Cells.Item(110, 1).Value = "Sri Lanka"
$ws.Cells.Item(110, 2).Value = 2511
$ws.Cells.Item(110, 3).Value = 57
$ws.Cells.Item(110, 4).Value = 1980
$ws.Cells.Item(110, 5).Value = 520
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 11

$ws.Cells.Item(112, 1).Value = "Mali"
$ws.Cells.Item(112, 2).Value = 2406
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 4).Value = 1683
$ws.Cells.Item(112, 5).Value = 602
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 121

$ws.Cells.Item(113, 1).Value = "Malaui"
$ws.Cells.Item(113, 2).Value = 2261
$ws.Cells.Item(113, 3).Value = 192
$ws.Cells.Item(113, 4).Value = 517
$ws.Cells.Item(113, 5).Value = 1711
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 2
$ws.Cells.Item(113, 8).Value = 33

$ws.Cells.Item(114, 1).Value = "Libano"
$ws.Cells.Item(114, 2).Value = 2168
$ws.Cells.Item(114, 3).Value = 86
$ws.Cells.Item(114, 4).Value = 1402
$ws.Cells.Item(114, 5).Value = 730
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 36

$ws.Cells.Item(130, 1).Value = "Suazilandia"
$ws.Cells.Item(130, 2).Value = 1311
$ws.Cells.Item(130, 3).Value = 54
$ws.Cells.Item(130, 4).Value = 656
$ws.Cells.Item(130, 5).Value = 637
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 18

$ws.Cells.Item(131, 1).Value = "Benin"
$ws.Cells.Item(131, 2).Value = 1285
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 333
$ws.Cells.Item(131, 5).Value = 929
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 23

$ws.Cells.Item(136, 1).Value = "Montenegro"
$ws.Cells.Item(136, 2).Value = 1164
$ws.Cells.Item(136, 3).Value = 145
$ws.Cells.Item(136, 4).Value = 325
$ws.Cells.Item(136, 5).Value = 816
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 4
$ws.Cells.Item(136, 8).Value = 23

$ws.Cells.Item(183, 1).Value = "Aruba"
$ws.Cells.Item(183, 2).Value = 105
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 99
$ws.Cells.Item(183, 5).Value = 3
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 3
